# Item 26 : Generics
# Duplicate slide 12 (Item 11, same "title slide" layout/animations) to
# create the new slide 13 (Item 26), then edit its text content in place.

$p = $ppt.ActivePresentation

$srcSlide = $p.Slides.Item(12)
$dupRange = $srcSlide.Duplicate()
$newSlide = $dupRange.Item(1)

# --- Title -------------------------------------------------------------
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Item 26"

# --- Subtitle ------------------------------------------------------------
$subtitleShape = $newSlide.Shapes.Item(2)
$subtitleTextRange = $subtitleShape.TextFrame.TextRange

# Paragraph 1
$subtitleTextRange.Paragraphs(1).Text = "How to declare a collection of objects?"

# Paragraph 2 - first run becomes the new sentence, then two more runs with
# their own formatting are appended after it.
$para2 = $subtitleTextRange.Paragraphs(2)
$run1 = $para2.Runs(1)
$run1.Text = "            Simply use raw collection:  " + [char]8220 + " private final Collection "

$run2 = $run1.InsertAfter("myList")

$run3 = $run2.InsertAfter("  = " + [char]8230 + ". " + [char]8221)

# Paragraph 6 (bullet 1)
$subtitleTextRange.Paragraphs(6).Text = "I want to avoid casting very object I read from the collection?"

# Paragraph 7 (bullet 2)
$subtitleTextRange.Paragraphs(7).Text = "I want to control what is inserted in my collection?"

# --- Rectangle callout ---------------------------------------------------
$rectShape = $newSlide.Shapes.Item(3)
$rectShape.Left = 3685658
$rectShape.Width = 4875565

$rectTextRange = $rectShape.TextFrame.TextRange
$rectTextRange.Text = "Generics (java 5) to the rescue"
